# Update shared strings used as "sample/placeholder" values in the
# example sheets so they match the latest table format.
$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item("ValidExampleData")
$wsVertical = $wb.Worksheets.Item("ValidExampleVerticalData")

# ValidExampleData: row 5 held "skip\n0000" (A5) and "all\r\n1111" (B5:I5)
$wsData.Range("A5").Value = "skip"
$wsData.Range("B5:I5").Value = "both"

# ValidExampleVerticalData: column D held the same two sample values
$wsVertical.Range("D2").Value = "skip"
$wsVertical.Range("D3:D10").Value = "both"

# Update the stored cursor/selection position on each sheet
$wsData.Range("D16").Select()
$wsVertical.Range("E20").Select()
